# Doing Updates for Financials
# Clears out the prior-period (column E) Balance Sheet figures on the PRNB
# sheet: most become "NA", a few collapse to 0 (matching the existing
# pattern used elsewhere in the sheet for blank/not-applicable figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRNB")

# --- Balance sheet: prior period (column E) values cleared to "NA" ---
$naRows = @(41, 43, 45, 46, 48, 52, 54, 57, 59, 60, 62, 66, 72, 76)
foreach ($r in $naRows) {
    $ws.Range("E$r").Value = "NA"
}

# --- Row 58 (Short/Current Long Term Debt): entire row collapses to 0 ---
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0

# --- Individual cells that collapse to numeric 0 instead of "NA" ---
$ws.Range("E61").Value = 0
$ws.Range("E70").Value = 0
